$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Low-level instruction" bullet: collapse the three runs (split apart by
#    a gramStart/gramEnd proofErr pair around "human, but") back into a
#    single run, which also drops the now-unneeded proofErr markers.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "understood by human, but is easier",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "understood by human, but is easier", 2)

# ---------------------------------------------------------------------------
# 2. Drop the stray "_GoBack" bookmark that used to sit in the
#    "Lexical analysis" paragraph. A new "_GoBack" bookmark will land inside
#    the newly-typed "Assembly" paragraph inserted below, matching where the
#    author's cursor ended up after the edit.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3. Insert the new "What Do Different Languages Do?" term list (Assembly,
#    C, Python, PHP, Perl, Java, Ruby, Go, Scratch, Lisp) right after the
#    "What Do Different Languages Do?" heading paragraph and before
#    "The Importance of C".
# ---------------------------------------------------------------------------
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "What Do Different Languages Do?`r") {
        $anchor = $d.Paragraphs($i)
        break
    }
}

$anchor.Range.InsertParagraphAfter()
$insertRange = $d.Paragraphs($anchor.Index + 1).Range

$blockXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Assembly (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>asm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>): a low-level programming lan</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/></w:rPr><w:t>guage in which there is a very strong (but often</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>not one-to- one) correspondence between the language and the architecture’s machine code</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>instructions</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• C: A general purpose, classic language that prevents most unintended operations</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Python: Created in 1991, an interpreted high-level programming language for general-purpose</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>programming</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• PHP: (Hypertext Preprocessor) a widely used open-source general-purpose scripting language</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>that is suited for web development and can be embedded into HTML</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Perl: Developed in the 1980s, used for automating systems, acting as a glue between different</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>computer systems and being a popular language for CGI programming.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Java: a general-purpose, concurrent, strong typed, class-based, object-oriented language.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Ruby: a dynamic, open-source language with a focus on simplicity and productivity</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Go Gopher: The Go language (developed 2009 at Google) is an open-source language used for</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>memory management and safety features. It is based on the Algol and C languages.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Scratch: A free visual programming language for children.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>• Lisp: Specified in 1958, is a family of high-level computer programming languages. Fortran is the</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>only older language (by one year).</w:t></w:r></w:p>'
$insertRange.InsertXML($blockXml)

# ---------------------------------------------------------------------------
# 4. The extra paragraphs above push the page breaks later in the doc, so
#    the lastRenderedPageBreak markers have to move too:
#      - off "The Corporate Object Revolution"
#      - onto "Jython"
#      - off "Briefly on the Huge Subject of Microsoft"
#      - onto "How Do You Pick a Programming Language?"
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "The Corporate Object Revolution`r") {
        $d.Paragraphs($i).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>The Corporate Object Revolution</w:t></w:r></w:p>')
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Jython: Version of Python designed to run inside of Java`r") {
        $d.Paragraphs($i).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Jython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Version of Python designed to run inside of Java</w:t></w:r></w:p>')
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Briefly on the Huge Subject of Microsoft`r") {
        $d.Paragraphs($i).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="160"/></w:pPr><w:r><w:t>Briefly on the Huge Subject of Microsoft</w:t></w:r></w:p>')
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "How Do You Pick a Programming Language?`r") {
        $d.Paragraphs($i).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="160"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>How Do You Pick a Programming Language?</w:t></w:r></w:p>')
        break
    }
}

Write-Output "edit complete"
